$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 9 -Day 28).Date
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 1250
$ws.Range("P2").Value = 417

$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 9 -Day 28).Date
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 1000
$ws.Range("P3").Value = 333

$ws.Range("D4").Value = (Get-Date -Year 2023 -Month 8 -Day 10).Date

$ws.Range("D5").Value = (Get-Date -Year 2023 -Month 8 -Day 10).Date
$ws.Range("J5").Value = 60

$ws.Range("D6").Value = (Get-Date -Year 2023 -Month 8 -Day 25).Date
$ws.Range("J6").Value = 100

$ws.Range("D7").Value = (Get-Date -Year 2023 -Month 8 -Day 22).Date
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("P7").Value = 833

$ws.Range("D8").Value = (Get-Date -Year 2023 -Month 9 -Day 7).Date
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2500
$ws.Range("P8").Value = 833

$ws.Range("D9").Value = (Get-Date -Year 2023 -Month 7 -Day 27).Date
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2500
$ws.Range("P9").Value = 833

$ws.Range("D10").Value = (Get-Date -Year 2023 -Month 8 -Day 28).Date
$ws.Range("J10").Value = 120

$ws.Range("D11").Value = (Get-Date -Year 2022 -Month 10 -Day 4).Date
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1250
$ws.Range("P11").Value = 417

$ws.Range("D12").Value = (Get-Date -Year 2022 -Month 10 -Day 4).Date
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 1000
$ws.Range("P12").Value = 333

$ws.Range("D13").Value = (Get-Date -Year 2023 -Month 8 -Day 23).Date
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2500
$ws.Range("P13").Value = 833

$ws.Range("D14").Value = (Get-Date -Year 2023 -Month 7 -Day 26).Date
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("P14").Value = 833

$ws.Range("D15").Value = (Get-Date -Year 2023 -Month 8 -Day 8).Date
$ws.Range("J15").Value = 80

$ws.Range("D16").Value = (Get-Date -Year 2023 -Month 8 -Day 8).Date
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 2000
$ws.Range("P16").Value = 667

$ws.Range("D17").Value = (Get-Date -Year 2023 -Month 8 -Day 7).Date
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2500
$ws.Range("P17").Value = 833

$ws.Range("D18").Value = (Get-Date -Year 2023 -Month 8 -Day 7).Date
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 2000
$ws.Range("P18").Value = 667

$ws.Range("D19").Value = (Get-Date -Year 2023 -Month 7 -Day 28).Date
$ws.Range("J19").Value = 70

$ws.Range("D20").Value = (Get-Date -Year 2022 -Month 10 -Day 12).Date
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 1200
$ws.Range("L20").Value = 1300
$ws.Range("M20").Value = 1250
$ws.Range("P20").Value = 417

$ws.Range("D21").Value = (Get-Date -Year 2022 -Month 10 -Day 12).Date
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 1000
$ws.Range("P21").Value = 333

$ws.Range("D22").Value = (Get-Date -Year 2023 -Month 8 -Day 11).Date
$ws.Range("J22").Value = 80

$ws.Range("D23").Value = (Get-Date -Year 2023 -Month 8 -Day 11).Date
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = 2000
$ws.Range("P23").Value = 667

$ws.Range("D24").Value = (Get-Date -Year 2023 -Month 9 -Day 6).Date
$ws.Range("J24").Value = 150

$ws.Range("D25").Value = (Get-Date -Year 2023 -Month 9 -Day 4).Date
$ws.Range("J25").Value = 100
